$wb = $excel.ActiveWorkbook

# Sheet "Total Hits"
$ws = $wb.Worksheets.Item("Total Hits")
$ws.Range("B2").Value = 1170
$ws.Range("D2").Value = "'40.18%"
$ws.Range("B3").Value = 2339
$ws.Range("D3").Value = "'40.16%"
$ws.Range("B4").Value = 3477
$ws.Range("D4").Value = "'39.80%"
$ws.Range("B5").Value = 4639
$ws.Range("D5").Value = "'39.83%"
$ws.Range("B6").Value = 5823
$ws.Range("D6").Value = "'39.99%"

# Sheet "Hits_entity"
$ws = $wb.Worksheets.Item("Hits_entity")
$ws.Range("B2").Value = 580
$ws.Range("D2").Value = "'41.11%"
$ws.Range("B3").Value = 1151
$ws.Range("D3").Value = "'40.79%"
$ws.Range("B4").Value = 1740
$ws.Range("D4").Value = "'41.11%"
$ws.Range("B5").Value = 2342
$ws.Range("D5").Value = "'41.50%"
$ws.Range("B6").Value = 2926
$ws.Range("D6").Value = "'41.47%"

# Sheet "Hits_numerical"
$ws = $wb.Worksheets.Item("Hits_numerical")
$ws.Range("B2").Value = 169
$ws.Range("D2").Value = "'25.80%"
$ws.Range("B3").Value = 334
$ws.Range("D3").Value = "'25.50%"
$ws.Range("B4").Value = 485
$ws.Range("D4").Value = "'24.68%"
$ws.Range("B5").Value = 629
$ws.Range("D5").Value = "'24.01%"
$ws.Range("B6").Value = 809
$ws.Range("D6").Value = "'24.70%"

# Sheet "Hits_boolean"
$ws = $wb.Worksheets.Item("Hits_boolean")
$ws.Range("B2").Value = 308
$ws.Range("D2").Value = "'53.75%"
$ws.Range("B3").Value = 631
$ws.Range("D3").Value = "'55.06%"
$ws.Range("B4").Value = 924
$ws.Range("D4").Value = "'53.75%"
$ws.Range("B5").Value = 1228
$ws.Range("D5").Value = "'53.58%"
$ws.Range("B6").Value = 1533
$ws.Range("D6").Value = "'53.51%"

# Sheet "Hits_date"
$ws = $wb.Worksheets.Item("Hits_date")
$ws.Range("B2").Value = 109
$ws.Range("D2").Value = "'41.13%"
$ws.Range("B3").Value = 217
$ws.Range("D3").Value = "'40.94%"
$ws.Range("B4").Value = 320
$ws.Range("D4").Value = "'40.25%"
$ws.Range("B5").Value = 428
$ws.Range("D5").Value = "'40.38%"
$ws.Range("B6").Value = 539
$ws.Range("D6").Value = "'40.68%"

# Sheet "Hits_string"
$ws = $wb.Worksheets.Item("Hits_string")
$ws.Range("B3").Value = 6
$ws.Range("D3").Value = "'37.50%"
$ws.Range("B4").Value = 8
$ws.Range("D4").Value = "'33.33%"
$ws.Range("B5").Value = 12
$ws.Range("D5").Value = "'37.50%"
$ws.Range("B6").Value = 16
$ws.Range("D6").Value = "'40.00%"
